$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Data": append new weekly WTREGEN observations (rows 104-109)
# ---------------------------------------------------------------
$wsData = $wb.Worksheets.Item("Data")

$newRows = @(
    @{ Row = 104; Date = 45189; Value = 640.296 },
    @{ Row = 105; Date = 45196; Value = 681.143 },
    @{ Row = 106; Date = 45203; Value = 668.009 },
    @{ Row = 107; Date = 45210; Value = 709.16 },
    @{ Row = 108; Date = 45217; Value = 759.878 },
    @{ Row = 109; Date = 45224; Value = 834.418 }
)

foreach ($r in $newRows) {
    $dateCell = $wsData.Cells.Item($r.Row, 1)
    $valCell = $wsData.Cells.Item($r.Row, 2)

    # Copy formatting from the row above (same date style used throughout column A)
    $wsData.Cells.Item($r.Row - 1, 1).Copy($dateCell)

    $dateCell.Value = $r.Date
    $valCell.Value = $r.Value
}

# ---------------------------------------------------------------
# Sheet "SeriesInfo": refresh metadata to match the new data pull
# ---------------------------------------------------------------
$wsInfo = $wb.Worksheets.Item("SeriesInfo")

# Leading apostrophe forces these ISO-date-looking strings to stay text,
# matching the original inlineStr cell type instead of becoming date serials.
$wsInfo.Range("B3").Formula = "'2023-10-27"
$wsInfo.Range("B4").Formula = "'2023-10-27"
$wsInfo.Range("B7").Formula = "'2023-10-25"
$wsInfo.Range("B14").Value = "2023-10-26 15:34:02-05"
$wsInfo.Range("B15").Value = 78
